$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.476.27"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.555.28"
$ws.Range("E3").Value = "  +4.88%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.85"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.84"
$ws.Range("E6").Value = "  +7.76%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "2.552.86"
$ws.Range("E9").Value = "  +4.80%  "
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.19"
$ws.Range("E14").Value = "  +7.24%  "
$ws.Range("D15").Value = "3.011.82"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").Value = "63.354.62"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "2.550.12"
$ws.Range("E18").Value = "  +4.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.61"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.34"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.97"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("E26").Value = "  +13.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.61"
$ws.Range("E27").Value = "  +5.22%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.49"
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("E30").Value = "  +12.89%  "
$ws.Range("D31").Value = "0.0₃0820"
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.86"
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "178.11"
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "424.11"
$ws.Range("E35").Value = "  +10.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.403"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.02"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.55"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "153.88"
$ws.Range("E43").Value = "  +6.18%  "
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.94"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.611"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0968"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0525"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").Value = "  +7.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.59"
$ws.Range("E50").Value = "  +3.98%  "
